$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 / 8 column B: cells removed ------------------------------------
$ws.Range("B6").ClearContents()
$ws.Range("B8").ClearContents()

# --- Prepare the new rows 29-30 / 33-35 so they carry the same look as
#     the existing "s0..s9" label rows (style index 1: light fill). Format
#     is copied first so the cell values below are written into cells
#     that already have their final formatting (keeps shared-string
#     insertion order aligned with the authoring order of the real edit).
$ws.Range("A19").Copy()
$ws.Range("A29:A30").PasteSpecial(-4122)
$ws.Range("A19").Copy()
$ws.Range("A33:A35").PasteSpecial(-4122)

# --- Write the new cell values in the same order the original author
#     introduced them (this drives the shared-strings insertion order).
$ws.Range("B35").Value = "R"
$ws.Range("B34").Value = "G"
$ws.Range("B33").Value = "B"
$ws.Range("A29").Value = "s10"
$ws.Range("A30").Value = "s11"
$ws.Range("B13").Value = "output"
$ws.Range("B7").Value = "last pointer new array"
$ws.Range("A33").Value = "d8"
$ws.Range("A34").Value = "d9"
$ws.Range("A35").Value = "d10"

# --- Row 32: single blank, explicitly-filled-none cell --------------------
$ws.Range("A32").Interior.ColorIndex = -4142

# --- Rows 36-47: blank cells using a Calibri-only font + explicit fill ----
foreach ($r in 36..47) {
    $cell = $ws.Range("A$r")
    $cell.Font.Name = "Calibri"
    $cell.Interior.ColorIndex = -4142
}

# --- Selection moves from G5 to B3 ----------------------------------------
$ws.Range("B3").Select() | Out-Null
